# Update the "Forecast Comparison" sheet:
#   - insert a new "Week_Start_Date" column right after "Week" (shifts
#     ASIN..is_holiday_week one column to the right, B..I -> C..J)
#   - fill the new column with each row's week-start date (as literal text)
#   - normalize the "Week" labels from W01..W09 to W1..W9 (no leading zero)
#   - re-enter is_holiday_week as a proper boolean (FALSE) instead of 0/1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new blank column at B; everything from B onward (ASIN, MyForecast,
# Amazon Mean/P70/P80/P90 Forecast, Product Title, is_holiday_week) shifts
# right by one, taking its data/formatting along with it.
$ws.Columns.Item(2).Insert()

# New header for the inserted column.
$ws.Range("B1").Value = "Week_Start_Date"

# Week labels (column A) without the leading zero, and week-start dates
# (column B), for each of the 16 data rows.
$weeks = @(
    @{Row=2;  Label="W1";  Date="2025-01-05"},
    @{Row=3;  Label="W2";  Date="2025-01-12"},
    @{Row=4;  Label="W3";  Date="2025-01-19"},
    @{Row=5;  Label="W4";  Date="2025-01-26"},
    @{Row=6;  Label="W5";  Date="2025-02-02"},
    @{Row=7;  Label="W6";  Date="2025-02-09"},
    @{Row=8;  Label="W7";  Date="2025-02-16"},
    @{Row=9;  Label="W8";  Date="2025-02-23"},
    @{Row=10; Label="W9";  Date="2025-03-02"},
    @{Row=11; Label="W10"; Date="2025-03-09"},
    @{Row=12; Label="W11"; Date="2025-03-16"},
    @{Row=13; Label="W12"; Date="2025-03-23"},
    @{Row=14; Label="W13"; Date="2025-03-30"},
    @{Row=15; Label="W14"; Date="2025-04-06"},
    @{Row=16; Label="W15"; Date="2025-04-13"},
    @{Row=17; Label="W16"; Date="2025-04-20"}
)

foreach ($week in $weeks) {
    $ws.Range("A$($week.Row)").Value = $week.Label
    # Leading apostrophe forces the ISO date string to be stored as plain
    # text instead of being auto-converted into a date serial number.
    $ws.Range("B$($week.Row)").Value = "'" + $week.Date
    # is_holiday_week (now column J after the insert) should be a real
    # boolean, not a 0/1 number.
    $ws.Range("J$($week.Row)").Value = $false
}
